$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B5 from a literal number to the text placeholder "{{totalParticip}}"
$ws.Range("B5").Value = "{{totalParticip}}"

# Update the active selection (view state) to C15
$ws.Range("C15").Select()
